# Update the "dSF" column (column F) values per the repulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -2
    4  = -4
    5  = -4
    6  = 3
    7  = -4
    8  = 1
    9  = -8
    10 = 5
    11 = -8
    12 = -1
    13 = 3
    14 = 4
    15 = -4
    16 = 6
    17 = 4
    18 = -3
    19 = -3
    21 = -4
    22 = -4
    24 = 4
    25 = -4
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
